# Update the Hardware-ID license sheet with the latest data.
#
# Original data (rows 2-4, data rows under the header):
#   Row2: S36SNWAH859775t | 2026-02-08 (46061)
#   Row3: S36SNWAH859775t | 2026-02-09 (46062)
#   Row4: S36SNWAH859775Z | 2050-02-10 (54829)
#   Row5: (empty trailing row)
#
# New data:
#   Row2: 002538B28101C704 | 2050-02-08 (54827)
#   Row3: S36SNWAH859775Z  | 2050-02-10 (54829)   (unchanged, shifted up)
#   Row4: (empty trailing row, shifted up)
#
# i.e. the stale duplicate "S36SNWAH859775t" row is removed, and the
# remaining first data row is refreshed with the new Hardware-ID reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete duplicate row (old row 3); rows below shift up.
$ws.Rows(3).Delete()

# Refresh the first data row with the newest Hardware-ID + date values.
$ws.Range("A2").Value = "002538B28101C704"
$ws.Range("B2").Value = 54827

# Move the active selection to match the saved workbook state.
$ws.Range("B16").Select()
